$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# D-column (Price) values are forced to text (NumberFormat "@" + Style reset)
# so numeric-looking strings (e.g. "1.0000", "0.9997") keep their exact
# textual representation instead of being coerced to a Double by Excel.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '25.964.09'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.76%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.746.23'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.22%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9997'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '248.71'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +5.18%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.0000'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5044'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -4.83%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2744'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.36%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06185'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.08%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.07276'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.32%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.743.59'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '15.20'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.6539'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.97%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '4.649'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.31%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '77.71'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.98%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '25.980.85'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('E19').Value = '  +0.45%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.000006837'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.99%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.967.70'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.18%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.423'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.89%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '8.722'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.13%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.397'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.98%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '136.82'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -2.05%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.505'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('E27').Value = '  -0.27%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.780'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.43%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '105.51'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.45%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '3.861'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +2.29%  '
$ws.Range('E31').Value = '  -1.14%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.635'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.07%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.04673'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('E34').Value = '  +0.40%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.9940'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.46%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.6199'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.86%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.751'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.43%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.01613'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.75%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.927'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.58%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.01%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '100.54'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.40%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.3939'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.33%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.7597'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.84%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '5.007'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.56%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.1149'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.31%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '6.310'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '55.65'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.83%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.05273'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.26%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '30.67'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.21%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '7.576'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.20%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.3436'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.39%  '
